# Update "想去人数" (F) / "最低票价" (G) figures to the latest scrape values.
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1): rows 3-12 ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 7556
$wsExhibit.Range("G3").Value = 79
$wsExhibit.Range("F4").Value = 285
$wsExhibit.Range("F5").Value = 18
$wsExhibit.Range("F6").Value = 459
$wsExhibit.Range("F7").Value = 4176
$wsExhibit.Range("F8").Value = 326
$wsExhibit.Range("F9").Value = 581
$wsExhibit.Range("F10").Value = 278
$wsExhibit.Range("F11").Value = 669
$wsExhibit.Range("F12").Value = 154

# --- Sheet "全部类型" (sheet4): corresponding rows ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 7556
$wsAll.Range("G4").Value = 79
$wsAll.Range("F6").Value = 285
$wsAll.Range("F7").Value = 18
$wsAll.Range("F8").Value = 459
$wsAll.Range("F9").Value = 4176
$wsAll.Range("F10").Value = 326
$wsAll.Range("F11").Value = 581
$wsAll.Range("F12").Value = 278
$wsAll.Range("F13").Value = 669
$wsAll.Range("F15").Value = 154
